$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a unique run of text anywhere in the document body.
# ---------------------------------------------------------------------------
function Replace-UniqueText($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Monte-Carlo run counts / random seed paragraph near the top of the doc.
# ---------------------------------------------------------------------------
Replace-UniqueText "For Retrofit (Alternative 1) 13000 Monte-Carlo simulations were run." "For Retrofit (Alternative 1) 12000 Monte-Carlo simulations were run."
Replace-UniqueText "For New Bridge (Alternative 2) 12000 Monte-Carlo simulations were run." "For New Bridge (Alternative 2) 9000 Monte-Carlo simulations were run."
Replace-UniqueText "The random number seed for these runs was 1622842854." "The random number seed for these runs was 1845552156."

# ---------------------------------------------------------------------------
# 2. Summary table (Tables.Item(1)) edits.
#    Row numbering (1-based, includes header row 1):
#      1  header (blank / Base / Retrofit / New Bridge)
#      2  Total Benefits ($)
#      3  (Lower Bound, Upper Bound) ($)          <- Total Benefits bound
#      4  Total Costs ($)
#      5  (Lower Bound, Upper Bound) ($)          <- Total Costs bound
#      6  Externalities ($)
#      7  (Lower Bound, Upper Bound) ($)          <- Externalities bound
#      8  Net with externalities ($)
#      9  (Lower Bound, Upper Bound) ($)          <- Net w/ ext bound
#      10 SIR with externalities(%)     -> BCR (%)
#      11 (Lower Bound, Upper Bound) ($)          <- SIR/BCR bound
#      12 IRR with externalities (%)    -> IRR (%)
#      13 (Lower Bound, Upper Bound) ($)          <- IRR bound
#      14 ROI with externalities (%)    -> ROI (%)
#      15 (Lower Bound, Upper Bound) ($)          <- ROI bound
#      16 Non-Disaster ROI with externalities (%) -> Non-Disaster ROI (%)
#      17 (Lower Bound, Upper Bound) ($)          <- Non-Disaster ROI bound
#      18 Net ($)                                  \
#      19 (Lower Bound, Upper Bound) ($)            |
#      20 SIR (%)                                   |
#      21 (Lower Bound, Upper Bound) ($)            |  duplicate block
#      22 IRR (%)                                   |  (no externalities)
#      23 (Lower Bound, Upper Bound) ($)            |  -- deleted --
#      24 ROI (%)                                   |
#      25 (Lower Bound, Upper Bound) ($)            |
#      26 Non-Disaster ROI (%)                      |
#      27 (Lower Bound, Upper Bound) ($)            /
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# -- numeric bound updates --------------------------------------------------
$t.Cell(3,3).Range.Text  = "(2,498,890; 5,081,529)"
$t.Cell(3,4).Range.Text  = "(6,432,061; 10,952,521)"

$t.Cell(5,3).Range.Text  = "(3,452,322; 4,288,263)"
$t.Cell(5,4).Range.Text  = "(7,712,053; 8,995,563)"

$t.Cell(9,3).Range.Text  = "(-1,411,533; 1,352,193)"
$t.Cell(9,4).Range.Text  = "(2,053,995; 6,788,917)"

# -- SIR with externalities -> BCR ------------------------------------------
$t.Cell(10,1).Range.Text = "BCR (%)"
$t.Cell(10,2).Range.Text = "No Valid BCR"

$t.Cell(11,2).Range.Text = "(No Valid BCR; No Valid BCR)"
$t.Cell(11,3).Range.Text = "(-0.36; 0.37)"
$t.Cell(11,4).Range.Text = "(0.26; 0.93)"

# -- IRR with externalities -> IRR -------------------------------------------
$t.Cell(12,1).Range.Text = "IRR (%)"

$t.Cell(13,3).Range.Text = "(0.83; 4.87)"
$t.Cell(13,4).Range.Text = "(4.31; 7.30)"

# -- ROI with externalities -> ROI -------------------------------------------
$t.Cell(14,1).Range.Text = "ROI (%)"

$t.Cell(15,3).Range.Text = "(1.28; 2.74)"
$t.Cell(15,4).Range.Text = "(2.48; 3.69)"

# -- Non-Disaster ROI with externalities -> Non-Disaster ROI -----------------
$t.Cell(16,1).Range.Text = "Non-Disaster ROI (%)"

# -- Remove the trailing duplicate (non-externality) block of 10 rows -------
for ($i = 0; $i -lt 10; $i++) {
    $t.Rows.Item(18).Delete()
}

# ---------------------------------------------------------------------------
# 3. Distribution-description paragraphs: drop the trailing ".00" from
#    whole-dollar figures.
# ---------------------------------------------------------------------------
Replace-UniqueText "Retrofit Indirect Loss Reduction: Gaussian distribution with standard deviation of 600000.00" "Retrofit Indirect Loss Reduction: Gaussian distribution with standard deviation of 600000"
Replace-UniqueText "Retrofit Response and Recovery: Gaussian distribution with standard deviation of 180000.00" "Retrofit Response and Recovery: Gaussian distribution with standard deviation of 180000"

Replace-UniqueText "Retrofit Indirect Cost: Triangular distribution with a min of 475000.00 and a max of 750000.00" "Retrofit Indirect Cost: Triangular distribution with a min of 475000 and a max of 750000"
Replace-UniqueText "Retrofit Direct Cost: Triangular distribution with a min of 2850000.00 and a max of 3840000.00" "Retrofit Direct Cost: Triangular distribution with a min of 2850000 and a max of 3840000"

Replace-UniqueText "New Bridge Indirect Loss Reduction: Gaussian distribution with standard deviation of 1050000.00" "New Bridge Indirect Loss Reduction: Gaussian distribution with standard deviation of 1050000"
Replace-UniqueText "New Bridge Response and Recovery: Gaussian distribution with standard deviation of 300000.00" "New Bridge Response and Recovery: Gaussian distribution with standard deviation of 300000"

Replace-UniqueText "Reduced Commute Time: Triangular distribution with a min of 70000.00 and a max of 115000.00" "Reduced Commute Time: Triangular distribution with a min of 70000 and a max of 115000"

Replace-UniqueText "New Bridge OMR: Rectangular distribution with a min of 21375.00 and a max of 30000.00" "New Bridge OMR: Rectangular distribution with a min of 21375 and a max of 30000"
Replace-UniqueText "Additional Roadwork Indirect Cost: Triangular distribution with a min of 114000.00 and a max of 144000.00" "Additional Roadwork Indirect Cost: Triangular distribution with a min of 114000 and a max of 144000"
Replace-UniqueText "Bridge Construction Indirect Cost: Triangular distribution with a min of 166250.00 and a max of 224000.00" "Bridge Construction Indirect Cost: Triangular distribution with a min of 166250 and a max of 224000"
Replace-UniqueText "Additional Roadwork Direct Cost: Triangular distribution with a min of 2375000.00 and a max of 3000000.00" "Additional Roadwork Direct Cost: Triangular distribution with a min of 2375000 and a max of 3000000"
Replace-UniqueText "Bridge Construction Direct Cost: Triangular distribution with a min of 4037500.00 and a max of 5440000.00" "Bridge Construction Direct Cost: Triangular distribution with a min of 4037500 and a max of 5440000"
Replace-UniqueText "Additional Roadwork OMR: Rectangular distribution with a min of 3500.00 and a max of 4250.00" "Additional Roadwork OMR: Rectangular distribution with a min of 3500 and a max of 4250"
